$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P5").Value2 = 'Pierogi'
$ws.Range("P6").Value2 = 'Bigos'
$ws.Range("P7").Value2 = 'Żurek'
$ws.Range("P8").Value2 = 'Placki Ziemniaczane'
$ws.Range("P9").Value2 = 'Golabki'
$ws.Range("P10").Value2 = 'Makowiec'
$ws.Range("P11").Value2 = 'Barszcz'
$ws.Range("P12").Value2 = 'Sernik'
$ws.Range("P13").Value2 = 'Kopytka'
$ws.Range("P14").Value2 = 'Kaczka'
$ws.Range("P15").Value2 = 'Sałatka Jarzynowa'
$ws.Range("P16").Value2 = 'Fasolka Po Bretonsku'
$ws.Range("P17").Value2 = 'Nalesniki'
$ws.Range("P18").Value2 = 'Sernik na Zimno'
$ws.Range("P19").Value2 = 'Zrazy'
$ws.Range("P20").Value2 = 'Zupa Krem z Pomidorow'
$ws.Range("P21").Value2 = 'Krupnik'
$ws.Range("P22").Value2 = 'Torty'
$ws.Range("P23").Value2 = 'Zurek z Biala Kielbasa'
$ws.Range("P24").Value2 = 'Smazony Serek'
$ws.Range("S5").Value2 = 'Tom Yum Goong'
$ws.Range("S6").Value2 = 'Green Curry'
$ws.Range("S7").Value2 = 'Massaman Curry'
$ws.Range("S8").Value2 = 'Larb'
$ws.Range("S9").Value2 = 'Khao Pad'
$ws.Range("S10").Value2 = 'Gaeng Daeng'
$ws.Range("S11").Value2 = ' Satay'
$ws.Range("S12").Value2 = 'Thai Spring Rolls'
$ws.Range("S13").Value2 = 'Khao Soi'
$ws.Range("S14").Value2 = 'Thai Iced Tea'
$ws.Range("S15").Value2 = 'Moo Pad Krapow '
$ws.Range("S16").Value2 = 'Thai Fish Cakes'
$ws.Range("S17").Value2 = 'Pineapple Fried Rice'
$ws.Range("S18").Value2 = 'Gaeng Som '
$ws.Range("S19").Value2 = 'Jok '
$ws.Range("S20").Value2 = 'Sticky Rice with Mango'
$ws.Range("S21").Value2 = 'Roti'
$ws.Range("S22").Value2 = 'som tum'
$ws.Range("S23").Value2 = 'Panang curry'
$ws.Range("W5").Value2 = 'Nasi Lemak'
$ws.Range("W6").Value2 = 'Char Kway Teow'
$ws.Range("W7").Value2 = 'Beef Rendang'
$ws.Range("W8").Value2 = 'Chicken Satay'
$ws.Range("W9").Value2 = 'Roti Canai'
$ws.Range("W10").Value2 = 'Hainanese Chicken Rice'
$ws.Range("W11").Value2 = 'Mee Goreng'
$ws.Range("W12").Value2 = 'Nasi Kerabu'
$ws.Range("W13").Value2 = 'Asam Pedas'
$ws.Range("W14").Value2 = 'Kuih Lapis'
$ws.Range("W15").Value2 = 'Rendang Ayam'
$ws.Range("W16").Value2 = 'Nasi Goreng'
$ws.Range("W17").Value2 = 'Soto Ayam'
$ws.Range("W18").Value2 = 'Curry Puff'
